# Insert a new data row for "Ají" (Terminal Hortofrutícola Agro Chillán) at
# row 37, shifting all existing rows from 37 downward by one. The workbook's
# data is sorted by date (descending), and this represents a new weekly
# record being added to the top of the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 37 - this shifts rows 37..128 down to
# 38..129 and automatically extends the sheet dimension to A1:R129.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record's values.
$ws.Range("A37").Value = 7
$ws.Range("B37").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C37").Value = "Ñuble"
$ws.Range("D37").Value = 44949
$ws.Range("E37").Value = 16
$ws.Range("F37").Value = 100112021
$ws.Range("G37").Value = "Ají"
$ws.Range("H37").Value = "Americana (o)"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 50
$ws.Range("K37").Value = 12000
$ws.Range("L37").Value = 12000
$ws.Range("M37").Value = 12000
$ws.Range("N37").Value = "`$/caja 15 kilos"
$ws.Range("O37").Value = "Región del Maule"
$ws.Range("P37").Value = 800
$ws.Range("Q37").Value = 15
$ws.Range("R37").Value = "Hortaliza"
